# Auto-generated COM-interop script applying the weekly crime-data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number + reporting week dates) ---
$ws.Range("A8").Value = "Volume 32   Number  30"
$ws.Range("C9").Value = "Report Covering the Week  7/21/2025  Through  7/27/2025"

# --- Cells changing between numeric and "N/A" text placeholders ---
# Donor cells (unaffected by this edit) supply the exact style + shared-string
# combination Excel uses for these placeholder cells; .Copy() replicates style
# + type + value exactly (plain .Value assignment would coerce "0" back to a number).
$ws.Range("D14").Copy($ws.Range("G14"))   # -> s=13, t=s, "0"
$ws.Range("E14").Copy($ws.Range("H14"))   # -> s=13, t=s, "***.*"
$ws.Range("D14").Copy($ws.Range("C15"))   # -> s=13, t=s, "0"
$ws.Range("D14").Copy($ws.Range("C27"))   # -> s=13, t=s, "0"
$ws.Range("D14").Copy($ws.Range("D31"))   # -> s=13, t=s, "0"
$ws.Range("E14").Copy($ws.Range("E31"))   # -> s=13, t=s, "***.*"
$ws.Range("J14").Copy($ws.Range("F31"))   # -> s=14 (numeric style)
$ws.Range("F31").Value = 1

# --- Plain numeric refreshes (style unchanged) ---
$values = @{
    "N15" = -10
    "C16" = 3
    "D16" = 5
    "E16" = -40
    "F16" = 21
    "G16" = 15
    "H16" = 40
    "I16" = 132
    "J16" = 138
    "K16" = -4.347826086956
    "L16" = 69.230769230769
    "M16" = 20
    "N16" = -81.06169296987
    "D17" = 8
    "E17" = -75
    "F17" = 22
    "G17" = 39
    "H17" = -43.589743589743
    "I17" = 197
    "J17" = 200
    "K17" = -1.5
    "L17" = 34.013605442176
    "M17" = 212.698412698413
    "N17" = -15.811965811965
    "C18" = 4
    "D18" = 7
    "E18" = -42.857142857142
    "G18" = 19
    "H18" = -36.842105263157
    "I18" = 119
    "J18" = 126
    "K18" = -5.555555555555
    "L18" = -9.160305343511
    "M18" = 83.076923076923
    "N18" = -74.463519313304
    "C19" = 8
    "D19" = 10
    "E19" = -20
    "F19" = 52
    "G19" = 56
    "H19" = -7.142857142857
    "I19" = 378
    "J19" = 417
    "K19" = -9.352517985611
    "L19" = -5.5
    "M19" = 44.274809160305
    "N19" = -37.31343283582
    "C20" = 5
    "D20" = 3
    "E20" = 66.666666666666
    "F20" = 11
    "G20" = 13
    "H20" = -15.384615384615
    "I20" = 60
    "J20" = 49
    "K20" = 22.448979591836
    "L20" = 36.363636363636
    "M20" = 100
    "N20" = -85.294117647058
    "C21" = 22
    "D21" = 33
    "E21" = -33.333333333333
    "F21" = 119
    "G21" = 143
    "H21" = -16.783216783216
    "I21" = 895
    "J21" = 937
    "K21" = -4.482390608324
    "L21" = 11.45703611457
    "M21" = 68.233082706766
    "N21" = -63.047068538398
    "D22" = 2
    "E22" = -50
    "F22" = 7
    "G22" = 7
    "I22" = 36
    "J22" = 29
    "K22" = 24.137931034482
    "L22" = 2.857142857142
    "M22" = -5.263157894736
    "D23" = 2
    "G23" = 4
    "H23" = -25
    "J23" = 26
    "K23" = -23.076923076923
    "C24" = 43
    "D24" = 50
    "E24" = -14
    "F24" = 149
    "G24" = 183
    "H24" = -18.579234972677
    "I24" = 962
    "J24" = 1252
    "K24" = -23.162939297124
    "L24" = -17.70744225834
    "M24" = 12.383177570093
    "D25" = 45
    "E25" = -28.888888888888
    "F25" = 118
    "G25" = 172
    "H25" = -31.395348837209
    "I25" = 746
    "J25" = 1134
    "K25" = -34.2151675485
    "L25" = -28.952380952381
    "C26" = 6
    "D26" = 10
    "E26" = -40
    "F26" = 30
    "G26" = 32
    "H26" = -6.25
    "I26" = 245
    "J26" = 269
    "K26" = -8.921933085501
    "L26" = 13.425925925925
    "M26" = 20.689655172413
    "D28" = 1
    "E28" = 0
    "F28" = 7
    "H28" = 0
    "I28" = 51
    "J28" = 32
    "K28" = 59.375
    "L28" = 70
    "H31" = -50
    "I31" = 7
    "K31" = -12.5
    "L31" = -41.666666666666
}
foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}

